# Rename the "Summary" sheet to "Description" and keep the
# Print_Titles defined name (and the sheet's own print-title rows)
# in sync with the new sheet name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Description"
$ws.PageSetup.PrintTitleRows = "Description!`$1:`$4"
